# Applies two data-maintenance edits to the "DC-Colos" sheet:
#  1. Rows 55-58 (ZRH, LYS, BOD, SKP) are rotated so that ZRH moves from
#     the top of the block to the bottom (LYS, BOD, SKP, ZRH).
#  2. The CTU (Chengdu, China) row (row 264) is removed entirely, which
#     shifts every following row up by one and shrinks the used range
#     from A1:H331 to A1:H330.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rotate rows 55-58 -------------------------------------------------
$firstRow = 55
$lastRow = 58
$numRows = $lastRow - $firstRow + 1
$numCols = 8

$blockData = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value2)
    }
    $blockData += ,$rowVals
}

# Rotate left by one: the old first row becomes the new last row.
$rotated = @()
for ($i = 1; $i -lt $numRows; $i++) {
    $rotated += ,($blockData[$i])
}
$rotated += ,($blockData[0])

for ($i = 0; $i -lt $numRows; $i++) {
    $r = $firstRow + $i
    $rowVals = $rotated[$i]
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# --- 2. Delete the CTU (Chengdu, China) row -------------------------------
$ws.Rows.Item(264).Delete()
